# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-CellPlain($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

Set-CellText "D2" "28.295.37"
Set-CellText "E2" "  -0.92%  "
Set-CellText "D3" "1.550.54"
Set-CellText "E3" "  -1.42%  "
Set-CellText "E4" "  -0.12%  "
Set-CellText "D5" "209.56"
Set-CellText "E5" "  -1.59%  "
Set-CellText "E6" "  -1.39%  "
Set-CellText "E7" "  -0.07%  "
Set-CellText "E8" "  -2.20%  "
Set-CellText "E9" "  -2.13%  "
Set-CellText "E10" "  -1.33%  "
Set-CellText "E11" "  +0.18%  "
Set-CellText "D12" "1.772.21"
Set-CellText "E12" "  -1.39%  "
Set-CellText "D13" "1.541.98"
Set-CellText "E13" "  -1.80%  "
Set-CellText "D14" "28.289.67"
Set-CellText "E14" "  -0.86%  "
Set-CellText "E15" "  -1.48%  "
Set-CellText "E16" "  -2.50%  "
Set-CellText "D17" "60.55"
Set-CellText "E17" "  -2.87%  "
Set-CellText "D18" "226.99"
Set-CellText "E18" "  -1.52%  "
Set-CellText "E19" "  -0.84%  "
Set-CellText "E20" "  -2.64%  "
Set-CellText "E21" "  -0.09%  "
Set-CellText "D22" "3.90"
Set-CellText "E22" "  +0.29%  "
Set-CellText "D23" "8.83"
Set-CellText "E23" "  -3.15%  "
Set-CellText "E24" "  -6.23%  "
Set-CellText "D25" "149.44"
Set-CellText "E25" "  -1.49%  "
Set-CellText "E26" "  -1.77%  "
Set-CellText "E27" "  -0.94%  "
Set-CellText "E29" "  -3.27%  "
Set-CellText "E30" "  -4.23%  "
Set-CellText "E31" "  -4.38%  "
Set-CellText "E32" "  -1.57%  "
Set-CellText "E33" "  -2.00%  "
Set-CellText "D34" "1.381.78"
Set-CellText "E34" "  -0.63%  "
Set-CellText "E35" "  +0.93%  "
Set-CellText "D36" "1.48"
Set-CellText "E36" "  -3.48%  "
Set-CellText "E37" "  -1.32%  "
Set-CellText "D38" "2.58"
Set-CellText "E38" "  -1.40%  "
Set-CellText "E39" "  -2.99%  "
Set-CellPlain "B40" "RenderToken"
Set-CellPlain "C40" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D40" "1.91"
Set-CellText "E40" "  +1.11%  "
Set-CellPlain "B41" "ImmutableX"
Set-CellPlain "C41" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-CellText "D41" "0.509"
Set-CellText "E41" "  -2.99%  "
Set-CellPlain "B42" "PaxDollar"
Set-CellPlain "C42" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-CellText "D42" "1.00"
Set-CellText "E42" "  -0.11%  "
Set-CellText "D43" "0.775"
Set-CellText "E43" "  -1.96%  "
Set-CellText "E44" "  -1.29%  "
Set-CellText "E45" "  -2.27%  "
Set-CellText "D46" "61.79"
Set-CellText "E46" "  -1.96%  "
Set-CellPlain "B47" "WEMIXToken"
Set-CellPlain "C47" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText "D47" "0.907"
Set-CellText "E47" "  -6.26%  "
Set-CellPlain "B48" "RocketPoolETH"
Set-CellPlain "C48" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-CellText "D48" "1.685.72"
Set-CellText "E48" "  -1.30%  "
Set-CellText "D49" "85.44"
Set-CellText "E49" "  -1.28%  "
Set-CellText "D50" "42.24"
Set-CellText "E50" "  +6.15%  "
Set-CellText "E51" "  +0.40%  "
